# "changement chili con carne"
# Update the interactive sandwich-quantity calculator:
#  - Le Classique (B6) : 22 -> 20
#  - Le BLT        (B10): 25 -> 30
#  - Le Thon       (B12): 16 -> 12
#  - Remove the now unused "calibre gros / calibre X-gros / difference"
#    helper table (rows 21-23), which also drops the three related
#    shared strings automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the editable quantities; dependent formulas recalc automatically.
$ws.Range("B6").Value  = 20
$ws.Range("B10").Value = 30
$ws.Range("B12").Value = 12

# Drop the leftover "calibre" helper rows at the bottom of the sheet.
$ws.Rows("21:23").Delete()

# Match the selection left in the saved file (cells B14:B15).
$ws.Range("B14:B15").Select() | Out-Null
